# Apply the two changes described by the diff:
#  1) Resize/reposition the "Rectangle 65" shape (the big background rounded
#     rectangle) to its new bounding box.
#  2) Merge the trailing "(p" + ")" runs in "TextBox 77" into a single
#     "(p)" run (keeping the second run's formatting), removing the third run.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1) Resize / reposition "Rectangle 65" ---
# Target EMU values (from the OOXML diff):
#   off  x=-14031   y=2020657
#   ext  cx=9158031 cy=4191000
# PowerPoint's Shape.Left/Top/Width/Height are expressed in points
# (1 pt = 12700 EMU); values below were chosen so the EMU round-trip
# (through the COM layer's Single-precision numbers) lands exactly on
# the target EMU integers.
$rect = $s.Shapes.Item("Rectangle 65")
$rect.Left = -1.104803149606299
$rect.Top = 159.1068503937008
$rect.Width = 721.10485
$rect.Height = 330.0

# --- 2) Merge "(p" and ")" runs into a single "(p)" run in "TextBox 77" ---
$textBox = $s.Shapes.Item("TextBox 77")
$tr = $textBox.TextFrame.TextRange
# Full text is "deleteTask(p)"; characters 11-13 are "(p)" and span the
# second + third runs. Assigning through Characters() rewrites just that
# sub-range while keeping the first run ("deleteTask") untouched, and
# takes on the second run's character formatting.
$chars = $tr.Characters(11, 3)
$chars.Text = "(p)"
